# Commit: Wed, May 13, 2020  2:04:27 AM
#
# 1) Slide 6's table used the custom "Table_0" style; switch it to the
#    built-in PowerPoint table style instead.
# 2) The deck's main (slide-master) theme swaps its colour scheme from the
#    imported "Integral" palette to the standard Office palette.

$p = $ppt.ActivePresentation

# --- 1. Table style on slide 6 -------------------------------------------
$s6 = $p.Slides.Item(6)
for ($i = 1; $i -le $s6.Shapes.Count; $i++) {
    $shp = $s6.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{F2A50139-B5B6-400C-803C-E50887C7E6EB}")
    }
}

# --- 2. Theme colour scheme ------------------------------------------------
$theme = $p.SlideMaster.Theme
$cs = $theme.ThemeColorScheme

$cs.Colors(1).RGB  = 0        # dk1      000000
$cs.Colors(2).RGB  = 16777215 # lt1      FFFFFF
$cs.Colors(3).RGB  = 6968388   # dk2      44546A
$cs.Colors(4).RGB  = 15132391  # lt2      E7E6E6
$cs.Colors(5).RGB  = 13998939  # accent1  5B9BD5
$cs.Colors(6).RGB  = 3243501   # accent2  ED7D31
$cs.Colors(7).RGB  = 10855845  # accent3  A5A5A5
$cs.Colors(8).RGB  = 49407     # accent4  FFC000
$cs.Colors(9).RGB  = 12874308  # accent5  4472C4
$cs.Colors(10).RGB = 4697456   # accent6  70AD47
$cs.Colors(11).RGB = 12673797  # hlink    0563C1
$cs.Colors(12).RGB = 7491477   # folHlink 954F72
